{"js": "// Increase the \"space before\" on the heading paragraphs (T\u00edtulo 1 / T\u00edtulo\n// 2 styles) and enlarge the page top/bottom margins plus header/footer\n// distances, per the commit's layout-breathing-room tweak.\n//\n// Word stores spacing/margins in twips (1/20 pt) in the OOXML, but the\n// Word JS object model exposes them in points, so we divide the target\n// twip values by 20 before assigning them.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/styleBuiltIn\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.styleBuiltIn === Word.Style.heading1) {\n    // w:before 80 -> 120 twips (4pt -> 6pt)\n    paragraph.spaceBefore = 6;\n  } else if (paragraph.styleBuiltIn === Word.Style.heading2) {\n    // w:before 80 -> 100 twips (4pt -> 5pt)\n    paragraph.spaceBefore = 5;\n  }\n}\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const pageSetup = section.pageSetup;\n  // w:pgMar twips -> points (pt = twips / 20)\n  pageSetup.topMargin = 1531 / 20; // 76.55pt (was 1418 -> 70.9pt)\n  pageSetup.bottomMargin = 1531 / 20; // 76.55pt (was 1418 -> 70.9pt)\n  pageSetup.headerDistance = 794 / 20; // 39.7pt (was 737 -> 36.85pt)\n  pageSetup.footerDistance = 680 / 20; // 34pt (was 624 -> 31.2pt)\n}\n\nawait context.sync();\n", "ps1": "# Increase the \"space before\" on the heading paragraphs (T\u00edtulo 1 / T\u00edtulo\n# 2 styles) and enlarge the page top/bottom margins plus header/footer\n# distances, per the commit's layout-breathing-room tweak.\n#\n# Word's OOXML stores spacing/margins in twips (1/20 pt), but the Word\n# object model's ParagraphFormat.SpaceBefore / PageSetup.*Margin /\n# PageSetup.*Distance properties are expressed in points, so the target\n# twip values below are divided by 20 before being assigned.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Heading 1\") {\n        # w:before 80 -> 120 twips (4pt -> 6pt)\n        $p.Format.SpaceBefore = 6\n    } elseif ($styleName -eq \"Heading 2\") {\n        # w:before 80 -> 100 twips (4pt -> 5pt)\n        $p.Format.SpaceBefore = 5\n    }\n}\n\nforeach ($sec in $d.Sections) {\n    $ps = $sec.PageSetup\n    # w:pgMar twips -> points (pt = twips / 20)\n    $ps.TopMargin = 1531 / 20      # 76.55pt (was 1418 -> 70.9pt)\n    $ps.BottomMargin = 1531 / 20   # 76.55pt (was 1418 -> 70.9pt)\n    $ps.HeaderDistance = 794 / 20  # 39.7pt (was 737 -> 36.85pt)\n    $ps.FooterDistance = 680 / 20  # 34pt (was 624 -> 31.2pt)\n}\n"}
